# Update the cryptos list data (Price and Volume(1h) columns) for each coin row.
# NumberFormat is forced to Text ("@") before assigning the value so that
# numeric-looking strings (e.g. "0.999", "6.97") are NOT auto-converted to
# the Number type by Excel -- these columns hold inline-string display text.
# ClearFormats() immediately afterwards removes the temporary "@" number
# format again so the cell keeps its original (default) style/format,
# exactly matching the source workbook which applies no explicit style to
# these cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $newValue
    $range.ClearFormats()
}

Set-TextValue "D2" '60.959.61'
Set-TextValue "E2" '  -0.07%  '
Set-TextValue "D3" '2.886.18'
Set-TextValue "E3" '  -1.33%  '
Set-TextValue "D4" '0.999'
Set-TextValue "E4" '  -0.01%  '
Set-TextValue "D5" '588.98'
Set-TextValue "E5" '  -0.33%  '
Set-TextValue "D6" '138.67'
Set-TextValue "E6" '  -5.74%  '
Set-TextValue "E7" '  +0.03%  '
Set-TextValue "D8" '0.493'
Set-TextValue "E8" '  -2.92%  '
Set-TextValue "D9" '7.01'
Set-TextValue "E9" '  +0.89%  '
Set-TextValue "D10" '0.138'
Set-TextValue "E10" '  -4.30%  '
Set-TextValue "D11" '0.427'
Set-TextValue "E11" '  -3.09%  '
Set-TextValue "E12" '  -3.81%  '
Set-TextValue "D13" '32.23'
Set-TextValue "E13" '  -4.37%  '
Set-TextValue "E14" '  -0.47%  '
Set-TextValue "D15" '3.364.91'
Set-TextValue "E15" '  -1.25%  '
Set-TextValue "D16" '60.836.21'
Set-TextValue "E16" '  -0.11%  '
Set-TextValue "D17" '2.864.70'
Set-TextValue "E17" '  -2.08%  '
Set-TextValue "D18" '6.48'
Set-TextValue "E18" '  -3.32%  '
Set-TextValue "D19" '425.24'
Set-TextValue "E19" '  -1.55%  '
Set-TextValue "D20" '13.18'
Set-TextValue "E20" '  -1.74%  '
Set-TextValue "D21" '0.654'
Set-TextValue "E21" '  -3.78%  '
Set-TextValue "D22" '6.92'
Set-TextValue "E22" '  -2.48%  '
Set-TextValue "D23" '79.85'
Set-TextValue "E23" '  -1.93%  '
Set-TextValue "D24" '10.38'
Set-TextValue "E24" '  -4.91%  '
Set-TextValue "D25" '0.999'
Set-TextValue "E25" '  -0.06%  '
Set-TextValue "D26" '2.06'
Set-TextValue "E26" '  -6.56%  '
Set-TextValue "D27" '11.37'
Set-TextValue "E27" '  -4.44%  '
Set-TextValue "E28" '  -3.23%  '
Set-TextValue "E29" '  -8.73%  '
Set-TextValue "D30" '6.61'
Set-TextValue "E30" '  -5.74%  '
Set-TextValue "E31" '  +0.03%  '
Set-TextValue "D32" '25.52'
Set-TextValue "E32" '  -4.36%  '
Set-TextValue "E33" '  -5.24%  '
Set-TextValue "D34" '0.0₃0838'
Set-TextValue "E34" '  -3.30%  '
Set-TextValue "E35" '  -4.38%  '
Set-TextValue "D36" '5.43'
Set-TextValue "E36" '  -3.99%  '
Set-TextValue "D37" '48.92'
Set-TextValue "E37" '  -2.24%  '
Set-TextValue "D38" '2.80'
Set-TextValue "E38" '  -6.91%  '
Set-TextValue "D39" '1.90'
Set-TextValue "E39" '  -4.44%  '
Set-TextValue "D40" '8.32'
Set-TextValue "E40" '  -2.87%  '
Set-TextValue "E41" '  -5.46%  '
Set-TextValue "D42" '0.266'
Set-TextValue "E42" '  -5.74%  '
Set-TextValue "D43" '38.08'
Set-TextValue "E43" '  -8.05%  '
Set-TextValue "D44" '2.659.40'
Set-TextValue "D45" '131.04'
Set-TextValue "E45" '  -2.31%  '
Set-TextValue "E46" '  -4.66%  '
Set-TextValue "D47" '352.20'
Set-TextValue "E47" '  -6.88%  '
Set-TextValue "E49" '  -4.11%  '
Set-TextValue "D50" '22.28'
Set-TextValue "E50" '  -6.43%  '
Set-TextValue "D51" '1.93'
Set-TextValue "E51" '  -4.14%  '
